$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: RJ45 connector changed to vertical variant (new part number)
$ws.Range("G5").Value = "538-95503-2881"
$ws.Range("C5").Value = "95503-2881"

# Row 8: ESD suppressor diode changed to new CMS part (SMAJ60A-TP by MCC)
$ws.Range("C8").Value = "SMAJ60A-TP"
$ws.Range("G8").Value = "833-SMAJ60A-TP"
$ws.Range("E8").Value = "MCC"
$ws.Range("H8").Value = 0.51

$ws.Range("D4").Select()
